# Rename the "Parameters" folder entry to "Configurations".
# The ProjectConfiguration worksheet has, on row 3:
#   A3 = "paramsFolder"   -> "configurationsFolder"
#   B3 = "Parameters/"    -> "Configurations/"
# The active selection also moved from C14 to A4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "Configurations/"
$ws.Range("A3").Value = "configurationsFolder"

$null = $ws.Range("A4").Select()
